# Auto-generated edit script: applies numeric corrections to the
# "currentAveragePrice*"/"LevePrice*"/"LeveProfit*" columns (H:N) across
# the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets, per the scheduled
# runner's refreshed market-board data.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 6786
$ws.Range("I132").Value = 7856.2354
$ws.Range("J132").Value = 2237.5
$ws.Range("K132").Value = 23568.7062
$ws.Range("L132").Value = 6712.5
$ws.Range("M132").Value = -21038.7062
$ws.Range("N132").Value = -11772.5
$ws.Range("H135").Value = 4289.857
$ws.Range("J135").Value = 7863.727
$ws.Range("L135").Value = 70773.54300000001
$ws.Range("N135").Value = -75843.54300000001
$ws.Range("H137").Value = 22733210
$ws.Range("I137").Value = 62501816
$ws.Range("K137").Value = 187505448
$ws.Range("M137").Value = -187502898
$ws.Range("H138").Value = 4465.0454
$ws.Range("J138").Value = 7530.8
$ws.Range("L138").Value = 22592.4
$ws.Range("N138").Value = -32872.4

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 6985.7144
$ws.Range("I28").Value = 6985.7144
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 6985.7144
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -6793.7144
$ws.Range("N28").ClearContents()
$ws.Range("H32").Value = 218406.95
$ws.Range("I32").Value = 485235.38
$ws.Range("K32").Value = 485235.38
$ws.Range("M32").Value = -484948.38
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H55").Value = 29682.334
$ws.Range("J55").Value = 34499.5
$ws.Range("L55").Value = 34499.5
$ws.Range("N55").Value = -35129.5
$ws.Range("H61").Value = 2707612.5
$ws.Range("I61").Value = 5500.9585
$ws.Range("K61").Value = 5500.9585
$ws.Range("M61").Value = -5288.9585
$ws.Range("H74").Value = 1550903.6
$ws.Range("I74").Value = 1920122
$ws.Range("K74").Value = 1920122
$ws.Range("M74").Value = -1919248
$ws.Range("H77").Value = 1550903.6
$ws.Range("I77").Value = 1920122
$ws.Range("K77").Value = 9600610
$ws.Range("M77").Value = -9596242
$ws.Range("H99").Value = 6985.7144
$ws.Range("I99").Value = 6985.7144
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 6985.7144
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -3990.7144
$ws.Range("N99").ClearContents()
$ws.Range("H110").Value = 1315.2727
$ws.Range("I110").Value = 989
$ws.Range("K110").Value = 989
$ws.Range("M110").Value = 1056
$ws.Range("H132").Value = 737727.4
$ws.Range("I132").Value = 835794.75
$ws.Range("K132").Value = 2507384.25
$ws.Range("M132").Value = -2504854.25
$ws.Range("H136").Value = 2707612.5
$ws.Range("I136").Value = 5500.9585
$ws.Range("K136").Value = 16502.8755
$ws.Range("M136").Value = -13952.8755

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 98.5
$ws.Range("I22").Value = 97
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 97
$ws.Range("L22").Value = 100
$ws.Range("M22").Value = 76
$ws.Range("N22").Value = -446
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H94").Value = 4471.5
$ws.Range("I94").Value = 5018.25
$ws.Range("K94").Value = 5018.25
$ws.Range("M94").Value = -4567.25
$ws.Range("H107").Value = 9259847
$ws.Range("J107").Value = 554.1667
$ws.Range("L107").Value = 554.1667
$ws.Range("N107").Value = -4394.1667
$ws.Range("H130").Value = 25000
$ws.Range("J130").Value = 25000
$ws.Range("L130").Value = 25000
$ws.Range("N130").Value = -35040
$ws.Range("H134").Value = 3630234.5
$ws.Range("I134").Value = 4849.0513
$ws.Range("K134").Value = 14547.1539
$ws.Range("M134").Value = -12012.1539

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2927007.2
$ws.Range("I31").Value = 4632429
$ws.Range("K31").Value = 4632429
$ws.Range("M31").Value = -4632134
$ws.Range("H34").Value = 2927007.2
$ws.Range("I34").Value = 4632429
$ws.Range("K34").Value = 4632429
$ws.Range("M34").Value = -4632227

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1548.0588
$ws.Range("I113").Value = 1328
$ws.Range("J113").Value = 1595.2142
$ws.Range("K113").Value = 3984
$ws.Range("L113").Value = 4785.642599999999
$ws.Range("M113").Value = -1814
$ws.Range("N113").Value = -9125.642599999999
$ws.Range("H129").Value = 8555517
$ws.Range("I129").Value = 961.2
$ws.Range("J129").Value = 13902114
$ws.Range("K129").Value = 2883.6
$ws.Range("L129").Value = 41706342
$ws.Range("M129").Value = 2116.4
$ws.Range("N129").Value = -41716342
$ws.Range("H131").Value = 5054553
$ws.Range("J131").Value = 7942052.5
$ws.Range("L131").Value = 23826157.5
$ws.Range("N131").Value = -23836237.5
$ws.Range("H137").Value = 6252.8
$ws.Range("J137").Value = 9388.444
$ws.Range("L137").Value = 28165.332
$ws.Range("N137").Value = -38365.33199999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 28000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 28000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H113").Value = 4038.8333
$ws.Range("I113").Value = 3729.6667
$ws.Range("J113").Value = 4554.1113
$ws.Range("K113").Value = 3729.6667
$ws.Range("L113").Value = 4554.1113
$ws.Range("M113").Value = -1559.6667
$ws.Range("N113").Value = -8894.1113
$ws.Range("H126").Value = 6385.1577
$ws.Range("I126").Value = 8327.083000000001
$ws.Range("J126").Value = 3056.1428
$ws.Range("K126").Value = 24981.249
$ws.Range("L126").Value = 9168.428400000001
$ws.Range("M126").Value = -22511.249
$ws.Range("N126").Value = -14108.4284
$ws.Range("H132").Value = 11172.976
$ws.Range("I132").Value = 9438.777
$ws.Range("K132").Value = 28316.331
$ws.Range("M132").Value = -25786.331

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3038.1482
$ws.Range("I40").Value = 2354.182
$ws.Range("J40").Value = 6047.6
$ws.Range("K40").Value = 2354.182
$ws.Range("L40").Value = 6047.6
$ws.Range("M40").Value = -2218.182
$ws.Range("N40").Value = -6319.6
$ws.Range("H55").Value = 819.9
$ws.Range("I55").Value = 839.9286
$ws.Range("J55").Value = 773.1667
$ws.Range("K55").Value = 839.9286
$ws.Range("L55").Value = 773.1667
$ws.Range("M55").Value = -666.9286
$ws.Range("N55").Value = -1119.1667
$ws.Range("H122").Value = 3614.9
$ws.Range("J122").Value = 4987.5
$ws.Range("L122").Value = 14962.5
$ws.Range("N122").Value = -19862.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 3367994.5
$ws.Range("I113").Value = 1153.0526
$ws.Range("K113").Value = 3459.1578
$ws.Range("M113").Value = -1289.1578
$ws.Range("H132").Value = 4904310.5
$ws.Range("I132").Value = 7248159.5
$ws.Range("K132").Value = 21744478.5
$ws.Range("M132").Value = -21741948.5

